$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 2-52 with the corrected naive forecaster values (dates shift back one period,
# and forecast ratios are recalculated). Row 53 of the old data is removed entirely.

$ws.Cells.Item(2, 1).Value = 39583
$ws.Cells.Item(2, 2).Value = 2008
$ws.Cells.Item(2, 3).ClearContents()
$ws.Cells.Item(2, 4).Value = 2009
$ws.Cells.Item(2, 5).Value = 1.003756253906252

$ws.Cells.Item(3, 1).Value = 39765
$ws.Cells.Item(3, 2).Value = 2008
$ws.Cells.Item(3, 3).ClearContents()
$ws.Cells.Item(3, 4).Value = 2009
$ws.Cells.Item(3, 5).Value = 0.8212989654785341

$ws.Cells.Item(4, 1).Value = 39948
$ws.Cells.Item(4, 2).Value = 2009
$ws.Cells.Item(4, 3).Value = 1.287693099940079
$ws.Cells.Item(4, 4).Value = 2010
$ws.Cells.Item(4, 5).Value = 1.224010362214401

$ws.Cells.Item(5, 1).Value = 40130
$ws.Cells.Item(5, 2).Value = 2009
$ws.Cells.Item(5, 3).Value = 1.218009596270675
$ws.Cells.Item(5, 4).Value = 2010
$ws.Cells.Item(5, 5).Value = 1.183007486132071

$ws.Cells.Item(6, 1).Value = 40310
$ws.Cells.Item(6, 2).Value = 2010
$ws.Cells.Item(6, 3).Value = 0.5167526861706184
$ws.Cells.Item(6, 4).Value = 2011
$ws.Cells.Item(6, 5).Value = 0.9718821796794952

$ws.Cells.Item(7, 1).Value = 40494
$ws.Cells.Item(7, 2).Value = 2010
$ws.Cells.Item(7, 3).Value = 0.5544720893820188
$ws.Cells.Item(7, 4).Value = 2011
$ws.Cells.Item(7, 5).Value = 1.015842920196763

$ws.Cells.Item(8, 1).Value = 40676
$ws.Cells.Item(8, 2).Value = 2011
$ws.Cells.Item(8, 3).Value = 1.187829657075357
$ws.Cells.Item(8, 4).Value = 2012
$ws.Cells.Item(8, 5).Value = 1.00065194548169

$ws.Cells.Item(9, 1).Value = 40862
$ws.Cells.Item(9, 2).Value = 2011
$ws.Cells.Item(9, 3).Value = 1.173294700162031
$ws.Cells.Item(9, 4).Value = 2012
$ws.Cells.Item(9, 5).Value = 0.9092565586104273

$ws.Cells.Item(10, 1).Value = 41044
$ws.Cells.Item(10, 2).Value = 2012
$ws.Cells.Item(10, 3).Value = 1.071158385438342
$ws.Cells.Item(10, 4).Value = 2013
$ws.Cells.Item(10, 5).Value = 0.912403143334517

$ws.Cells.Item(11, 1).Value = 41228
$ws.Cells.Item(11, 2).Value = 2012
$ws.Cells.Item(11, 3).Value = 1.180518841971723
$ws.Cells.Item(11, 4).Value = 2013
$ws.Cells.Item(11, 5).Value = 1.236730309040235

$ws.Cells.Item(12, 1).Value = 41409
$ws.Cells.Item(12, 2).Value = 2013
$ws.Cells.Item(12, 3).Value = 0.9553801317191413
$ws.Cells.Item(12, 4).Value = 2014
$ws.Cells.Item(12, 5).Value = 1.066801818459595

$ws.Cells.Item(13, 1).Value = 41592
$ws.Cells.Item(13, 2).Value = 2013
$ws.Cells.Item(13, 3).Value = 0.9276272455014611
$ws.Cells.Item(13, 4).Value = 2014
$ws.Cells.Item(13, 5).Value = 1.029015928490629

$ws.Cells.Item(14, 1).Value = 41774
$ws.Cells.Item(14, 2).Value = 2014
$ws.Cells.Item(14, 3).Value = 1.190496724073231
$ws.Cells.Item(14, 4).Value = 2015
$ws.Cells.Item(14, 5).Value = 1.154811676806311

$ws.Cells.Item(15, 1).Value = 41957
$ws.Cells.Item(15, 2).Value = 2014
$ws.Cells.Item(15, 3).Value = 1.265990289415564
$ws.Cells.Item(15, 4).Value = 2015
$ws.Cells.Item(15, 5).Value = 1.358148715145191

$ws.Cells.Item(16, 1).Value = 42137
$ws.Cells.Item(16, 2).Value = 2015
$ws.Cells.Item(16, 3).Value = 1.5464392869869
$ws.Cells.Item(16, 4).Value = 2016
$ws.Cells.Item(16, 5).Value = 1.247870081683522

$ws.Cells.Item(17, 1).Value = 42321
$ws.Cells.Item(17, 2).Value = 2015
$ws.Cells.Item(17, 3).Value = 1.642047742738506
$ws.Cells.Item(17, 4).Value = 2016
$ws.Cells.Item(17, 5).Value = 1.528208222695326

$ws.Cells.Item(18, 1).Value = 42503
$ws.Cells.Item(18, 2).Value = 2016
$ws.Cells.Item(18, 3).Value = 1.701952652941463
$ws.Cells.Item(18, 4).Value = 2017
$ws.Cells.Item(18, 5).Value = 1.637918813512695

$ws.Cells.Item(19, 1).Value = 42689
$ws.Cells.Item(19, 2).Value = 2016
$ws.Cells.Item(19, 3).Value = 1.66194179127146
$ws.Cells.Item(19, 4).Value = 2017
$ws.Cells.Item(19, 5).Value = 1.634928000057778

$ws.Cells.Item(20, 1).Value = 42867
$ws.Cells.Item(20, 2).Value = 2017
$ws.Cells.Item(20, 3).Value = 1.580693894992691
$ws.Cells.Item(20, 4).Value = 2018
$ws.Cells.Item(20, 5).Value = 1.610567777412109

$ws.Cells.Item(21, 1).Value = 43053
$ws.Cells.Item(21, 2).Value = 2017
$ws.Cells.Item(21, 3).Value = 1.609733807897773
$ws.Cells.Item(21, 4).Value = 2018
$ws.Cells.Item(21, 5).Value = 1.67176973076042

$ws.Cells.Item(22, 1).Value = 43145
$ws.Cells.Item(22, 2).Value = 2018
$ws.Cells.Item(22, 3).Value = 1.651937828695615
$ws.Cells.Item(22, 4).Value = 2019
$ws.Cells.Item(22, 5).Value = 1.63821551487775

$ws.Cells.Item(23, 1).Value = 43235
$ws.Cells.Item(23, 2).Value = 2018
$ws.Cells.Item(23, 3).Value = 1.646565058924154
$ws.Cells.Item(23, 4).Value = 2019
$ws.Cells.Item(23, 5).Value = 1.636439239090515

$ws.Cells.Item(24, 1).Value = 43326
$ws.Cells.Item(24, 2).Value = 2018
$ws.Cells.Item(24, 3).Value = 1.642460763882414
$ws.Cells.Item(24, 4).Value = 2019
$ws.Cells.Item(24, 5).Value = 1.62621273827539

$ws.Cells.Item(25, 1).Value = 43418
$ws.Cells.Item(25, 2).Value = 2018
$ws.Cells.Item(25, 3).Value = 1.641178243814534
$ws.Cells.Item(25, 4).Value = 2019
$ws.Cells.Item(25, 5).Value = 1.603287858019664

$ws.Cells.Item(26, 1).Value = 43510
$ws.Cells.Item(26, 2).Value = 2019
$ws.Cells.Item(26, 3).Value = 1.504616869537312
$ws.Cells.Item(26, 4).Value = 2020
$ws.Cells.Item(26, 5).Value = 1.599505522959732

$ws.Cells.Item(27, 1).Value = 43600
$ws.Cells.Item(27, 2).Value = 2019
$ws.Cells.Item(27, 3).Value = 1.619750436871126
$ws.Cells.Item(27, 4).Value = 2020
$ws.Cells.Item(27, 5).Value = 1.669486277487398

$ws.Cells.Item(28, 1).Value = 43691
$ws.Cells.Item(28, 2).Value = 2019
$ws.Cells.Item(28, 3).Value = 1.242963308065193
$ws.Cells.Item(28, 4).Value = 2020
$ws.Cells.Item(28, 5).Value = 1.082447181878954

$ws.Cells.Item(29, 1).Value = 43783
$ws.Cells.Item(29, 2).Value = 2019
$ws.Cells.Item(29, 3).Value = 1.183163144818633
$ws.Cells.Item(29, 4).Value = 2020
$ws.Cells.Item(29, 5).Value = 0.8408455317168162

$ws.Cells.Item(30, 1).Value = 43875
$ws.Cells.Item(30, 2).Value = 2020
$ws.Cells.Item(30, 3).Value = 0.6216637650511503
$ws.Cells.Item(30, 4).Value = 2021
$ws.Cells.Item(30, 5).Value = 1.126729649114599

$ws.Cells.Item(31, 1).Value = 43966
$ws.Cells.Item(31, 2).Value = 2020
$ws.Cells.Item(31, 3).Value = 0.2954364073068261
$ws.Cells.Item(31, 4).Value = 2021
$ws.Cells.Item(31, 5).Value = 0.8326407735962826

$ws.Cells.Item(32, 1).Value = 44068
$ws.Cells.Item(32, 2).Value = 2020
$ws.Cells.Item(32, 3).Value = -3.662861831460751
$ws.Cells.Item(32, 4).Value = 2021
$ws.Cells.Item(32, 5).Value = -1.110565553434917

$ws.Cells.Item(33, 1).Value = 44159
$ws.Cells.Item(33, 2).Value = 2020
$ws.Cells.Item(33, 3).Value = -3.662861831460751
$ws.Cells.Item(33, 4).Value = 2021
$ws.Cells.Item(33, 5).Value = -1.875058665585216

$ws.Cells.Item(34, 1).Value = 44251
$ws.Cells.Item(34, 2).Value = 2021
$ws.Cells.Item(34, 3).Value = -3.604628722764358
$ws.Cells.Item(34, 4).Value = 2022
$ws.Cells.Item(34, 5).Value = -2.856219939917704

$ws.Cells.Item(35, 1).Value = 44341
$ws.Cells.Item(35, 2).Value = 2021
$ws.Cells.Item(35, 3).Value = -1.564297238929013
$ws.Cells.Item(35, 4).Value = 2022
$ws.Cells.Item(35, 5).Value = 0.1548119563699935

$ws.Cells.Item(36, 1).Value = 44432
$ws.Cells.Item(36, 2).Value = 2021
$ws.Cells.Item(36, 3).Value = 0.1010915562932313
$ws.Cells.Item(36, 4).Value = 2022
$ws.Cells.Item(36, 5).Value = 6.64637963114707

$ws.Cells.Item(37, 1).Value = 44525
$ws.Cells.Item(37, 2).Value = 2021
$ws.Cells.Item(37, 3).Value = 0.1010915562932313
$ws.Cells.Item(37, 4).Value = 2022
$ws.Cells.Item(37, 5).Value = 5.03478667886097

$ws.Cells.Item(38, 1).Value = 44617
$ws.Cells.Item(38, 2).Value = 2022
$ws.Cells.Item(38, 3).Value = 5.220550987750228
$ws.Cells.Item(38, 4).Value = 2023
$ws.Cells.Item(38, 5).Value = 1.043506288584606

$ws.Cells.Item(39, 1).Value = 44706
$ws.Cells.Item(39, 2).Value = 2022
$ws.Cells.Item(39, 3).Value = 5.937304773291885
$ws.Cells.Item(39, 4).Value = 2023
$ws.Cells.Item(39, 5).Value = 2.150399152794202

$ws.Cells.Item(40, 1).Value = 44798
$ws.Cells.Item(40, 2).Value = 2022
$ws.Cells.Item(40, 3).Value = 5.793673192389748
$ws.Cells.Item(40, 4).Value = 2023
$ws.Cells.Item(40, 5).Value = 1.728278600643907

$ws.Cells.Item(41, 1).Value = 44890
$ws.Cells.Item(41, 2).Value = 2022
$ws.Cells.Item(41, 3).Value = 5.793673192389748
$ws.Cells.Item(41, 4).Value = 2023
$ws.Cells.Item(41, 5).Value = 2.399708479013141

$ws.Cells.Item(42, 1).Value = 44981
$ws.Cells.Item(42, 2).Value = 2023
$ws.Cells.Item(42, 3).Value = 0.156542203858212
$ws.Cells.Item(42, 4).Value = 2024
$ws.Cells.Item(42, 5).Value = 2.138492443986739

$ws.Cells.Item(43, 1).Value = 45071
$ws.Cells.Item(43, 2).Value = 2023
$ws.Cells.Item(43, 3).Value = -0.2621830498131694
$ws.Cells.Item(43, 4).Value = 2024
$ws.Cells.Item(43, 5).Value = 1.878976297039481

$ws.Cells.Item(44, 1).Value = 45163
$ws.Cells.Item(44, 2).Value = 2023
$ws.Cells.Item(44, 3).Value = -0.3788601787194756
$ws.Cells.Item(44, 4).Value = 2024
$ws.Cells.Item(44, 5).Value = 2.181307424743695

$ws.Cells.Item(45, 1).Value = 45254
$ws.Cells.Item(45, 2).Value = 2023
$ws.Cells.Item(45, 3).Value = -0.3788601787194756
$ws.Cells.Item(45, 4).Value = 2024
$ws.Cells.Item(45, 5).Value = 0.8520283695166997

$ws.Cells.Item(46, 1).Value = 45345
$ws.Cells.Item(46, 2).Value = 2024
$ws.Cells.Item(46, 3).Value = 0.005756553697899847
$ws.Cells.Item(46, 4).Value = 2025
$ws.Cells.Item(46, 5).Value = -0.1096192596443557

$ws.Cells.Item(47, 1).Value = 45436
$ws.Cells.Item(47, 2).Value = 2024
$ws.Cells.Item(47, 3).Value = 0.0512320434504332
$ws.Cells.Item(47, 4).Value = 2025
$ws.Cells.Item(47, 5).Value = 0.2660756331863467

$ws.Cells.Item(48, 1).Value = 45534
$ws.Cells.Item(48, 2).Value = 2024
$ws.Cells.Item(48, 3).Value = 0.05771202657300911
$ws.Cells.Item(48, 4).Value = 2025
$ws.Cells.Item(48, 5).Value = 0.2051249733294291

$ws.Cells.Item(49, 1).Value = 45618
$ws.Cells.Item(49, 2).Value = 2024
$ws.Cells.Item(49, 3).Value = 0.05771202657300911
$ws.Cells.Item(49, 4).Value = 2025
$ws.Cells.Item(49, 5).Value = 0.299857156820571

$ws.Cells.Item(50, 1).Value = 45713
$ws.Cells.Item(50, 2).Value = 2025
$ws.Cells.Item(50, 3).Value = 0.4540776569412763
$ws.Cells.Item(50, 4).Value = 2026
$ws.Cells.Item(50, 5).Value = -0.1815195499670796

$ws.Cells.Item(51, 1).Value = 45800
$ws.Cells.Item(51, 2).Value = 2025
$ws.Cells.Item(51, 3).Value = 0.5998844096825495
$ws.Cells.Item(51, 4).Value = 2026
$ws.Cells.Item(51, 5).Value = 0.1733734969819434

$ws.Cells.Item(52, 1).Value = 45891
$ws.Cells.Item(52, 2).Value = 2025
$ws.Cells.Item(52, 3).Value = 0.6062046309774693
$ws.Cells.Item(52, 4).Value = 2026
$ws.Cells.Item(52, 5).Value = 0.3877310837361314

# Remove the now-extra last row (row 53) so the sheet dimension matches A1:E52
$ws.Rows("53:53").Delete()
